# Add data for 2021-11-01: rename sheet to "Through 2021-10-24",
# update the October month row (row 12) and Total row (row 13)
# with updated no_arrest_made counts and recomputed arrest_rate values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename worksheet/tab title to reflect the new "through" date.
$ws.Name = "Through 2021-10-24"

# Update the label in A12.
$ws.Range("A12").Value = "October (through 10-24)"

# --- Row 12: October month totals ---
$ws.Range("C12").Value = 19
$ws.Range("D12").Value = 0.09520000000000001

$ws.Range("F12").Value = 37
$ws.Range("G12").Value = 0.075

$ws.Range("I12").Value = 47
$ws.Range("J12").Value = 0.1296

$ws.Range("O12").Value = 37
$ws.Range("P12").Value = 0.09760000000000001

$ws.Range("R12").Value = 120
$ws.Range("U12").Value = 152

# --- Row 13: Total row ---
$ws.Range("C13").Value = 215
$ws.Range("D13").Value = 0.1296

$ws.Range("F13").Value = 420
$ws.Range("G13").Value = 0.1045

$ws.Range("I13").Value = 624
$ws.Range("J13").Value = 0.0837

$ws.Range("O13").Value = 416
$ws.Range("P13").Value = 0.1015

$ws.Range("R13").Value = 968
$ws.Range("S13").Value = 0.0519

$ws.Range("U13").Value = 1317
$ws.Range("V13").Value = 0.0586
